# Auto-generated: append rows 1133-1183 to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(1133, 1).Value = "Buying Opportunity"
$ws.Cells.Item(1133, 2).Value = "support Zone"
$ws.Cells.Item(1133, 3).Value = "long buildup"
$ws.Cells.Item(1133, 4).Value = "Short buildup"
$ws.Cells.Item(1133, 5).Value = "FII ENTERING"

$ws.Cells.Item(1134, 1).Value = "AVALON"
$ws.Cells.Item(1134, 2).Value = "3IINFOLTD"
$ws.Cells.Item(1134, 3).Value = "MCX"
$ws.Cells.Item(1134, 4).Value = "BANDHANBNK"
$ws.Cells.Item(1134, 5).Value = "HDFCAMC"
$ws.Cells.Item(1134, 6).Value = 555.55
$ws.Cells.Item(1134, 7).Value = 36.45
$ws.Cells.Item(1134, 8).Value = 3941.6
$ws.Cells.Item(1134, 9).Value = 203.45
$ws.Cells.Item(1134, 10).Value = 4113.6

$ws.Cells.Item(1135, 1).Value = "COUNCODOS"
$ws.Cells.Item(1135, 2).Value = "ABCAPITAL"
$ws.Cells.Item(1135, 3).Value = "MUTHOOTFIN"
$ws.Cells.Item(1135, 4).Value = "COALINDIA"
$ws.Cells.Item(1135, 5).Value = "NAUKRI"
$ws.Cells.Item(1135, 6).Value = 6.17
$ws.Cells.Item(1135, 7).Value = 237.52
$ws.Cells.Item(1135, 8).Value = 1780.8
$ws.Cells.Item(1135, 9).Value = 469.25
$ws.Cells.Item(1135, 10).Value = 6648.6

$ws.Cells.Item(1136, 1).Value = "EMSLIMITED"
$ws.Cells.Item(1136, 2).Value = "ADVENZYMES"
$ws.Cells.Item(1136, 6).Value = 590.1
$ws.Cells.Item(1136, 7).Value = 380.5

$ws.Cells.Item(1137, 1).Value = "GULFOILLUB"
$ws.Cells.Item(1137, 2).Value = "AGROPHOS"
$ws.Cells.Item(1137, 6).Value = 1150.45
$ws.Cells.Item(1137, 7).Value = 46.49

$ws.Cells.Item(1138, 1).Value = "HDFCAMC"
$ws.Cells.Item(1138, 2).Value = "ANANDRATHI"
$ws.Cells.Item(1138, 6).Value = 4113.6
$ws.Cells.Item(1138, 7).Value = 3880.65

$ws.Cells.Item(1139, 1).Value = "INFOMEDIA"
$ws.Cells.Item(1139, 2).Value = "APOLLO"
$ws.Cells.Item(1139, 6).Value = 8.039999999999999
$ws.Cells.Item(1139, 7).Value = 108.69

$ws.Cells.Item(1140, 1).Value = "JBMA"
$ws.Cells.Item(1140, 2).Value = "ARIES"
$ws.Cells.Item(1140, 6).Value = 2148
$ws.Cells.Item(1140, 7).Value = 257.99

$ws.Cells.Item(1141, 1).Value = "KPIL"
$ws.Cells.Item(1141, 2).Value = "ASAL"
$ws.Cells.Item(1141, 6).Value = 1184.3
$ws.Cells.Item(1141, 7).Value = 1005.8

$ws.Cells.Item(1142, 1).Value = "NDGL"
$ws.Cells.Item(1142, 2).Value = "ASKAUTOLTD"
$ws.Cells.Item(1142, 6).Value = 3229.3
$ws.Cells.Item(1142, 7).Value = 352.25

$ws.Cells.Item(1143, 1).Value = "NEULANDLAB"
$ws.Cells.Item(1143, 2).Value = "BHAGYANGR"
$ws.Cells.Item(1143, 6).Value = 7342.05
$ws.Cells.Item(1143, 7).Value = 107.22

$ws.Cells.Item(1144, 1).Value = "NRBBEARING"
$ws.Cells.Item(1144, 2).Value = "BHARATRAS"
$ws.Cells.Item(1144, 6).Value = 349.3
$ws.Cells.Item(1144, 7).Value = 11695.7

$ws.Cells.Item(1145, 1).Value = "PGIL"
$ws.Cells.Item(1145, 2).Value = "BLS"
$ws.Cells.Item(1145, 6).Value = 756.45
$ws.Cells.Item(1145, 7).Value = 352.35

$ws.Cells.Item(1146, 1).Value = "PIXTRANS"
$ws.Cells.Item(1146, 2).Value = "BLUEJET"
$ws.Cells.Item(1146, 6).Value = 1489.95
$ws.Cells.Item(1146, 7).Value = 407.7

$ws.Cells.Item(1147, 1).Value = "PLASTIBLEN"
$ws.Cells.Item(1147, 2).Value = "COROMANDEL"
$ws.Cells.Item(1147, 6).Value = 275.05
$ws.Cells.Item(1147, 7).Value = 1525.35

$ws.Cells.Item(1148, 1).Value = "PNBGILTS"
$ws.Cells.Item(1148, 2).Value = "DALBHARAT"
$ws.Cells.Item(1148, 6).Value = 135.27
$ws.Cells.Item(1148, 7).Value = 1816.6

$ws.Cells.Item(1149, 1).Value = "PPAP"
$ws.Cells.Item(1149, 2).Value = "DCAL"
$ws.Cells.Item(1149, 6).Value = 221.62
$ws.Cells.Item(1149, 7).Value = 171.93

$ws.Cells.Item(1150, 1).Value = "QUESS"
$ws.Cells.Item(1150, 2).Value = "DEEPAKNTR"
$ws.Cells.Item(1150, 6).Value = 608.45
$ws.Cells.Item(1150, 7).Value = 2459.75

$ws.Cells.Item(1151, 1).Value = "RITCO"
$ws.Cells.Item(1151, 2).Value = "GNFC"
$ws.Cells.Item(1151, 6).Value = 293.69
$ws.Cells.Item(1151, 7).Value = 694.65

$ws.Cells.Item(1152, 1).Value = "RRKABEL"
$ws.Cells.Item(1152, 2).Value = "GSFC"
$ws.Cells.Item(1152, 6).Value = 1761.5
$ws.Cells.Item(1152, 7).Value = 238.04

$ws.Cells.Item(1153, 2).Value = "GSLSU"
$ws.Cells.Item(1153, 7).Value = 218.96

$ws.Cells.Item(1154, 2).Value = "GSPL"
$ws.Cells.Item(1154, 7).Value = 301.25

$ws.Cells.Item(1155, 2).Value = "GUJALKALI"
$ws.Cells.Item(1155, 7).Value = 787.4

$ws.Cells.Item(1156, 2).Value = "HINDPETRO"
$ws.Cells.Item(1156, 7).Value = 334.65

$ws.Cells.Item(1157, 2).Value = "INDOAMIN"
$ws.Cells.Item(1157, 7).Value = 125.22

$ws.Cells.Item(1158, 2).Value = "JBCHEPHARM"
$ws.Cells.Item(1158, 7).Value = 1744.4

$ws.Cells.Item(1159, 2).Value = "KTKBANK"
$ws.Cells.Item(1159, 7).Value = 226.56

$ws.Cells.Item(1160, 2).Value = "LAOPALA"
$ws.Cells.Item(1160, 7).Value = 328.35

$ws.Cells.Item(1161, 2).Value = "LICI"
$ws.Cells.Item(1161, 7).Value = 1007.65

$ws.Cells.Item(1162, 2).Value = "LINDEINDIA"
$ws.Cells.Item(1162, 7).Value = 8267.450000000001

$ws.Cells.Item(1163, 2).Value = "LTFOODS"
$ws.Cells.Item(1163, 7).Value = 263.43

$ws.Cells.Item(1164, 2).Value = "LUXIND"
$ws.Cells.Item(1164, 7).Value = 1485.7

$ws.Cells.Item(1165, 2).Value = "MADRASFERT"
$ws.Cells.Item(1165, 7).Value = 106.4

$ws.Cells.Item(1166, 2).Value = "MANAKSIA"
$ws.Cells.Item(1166, 7).Value = 97.72

$ws.Cells.Item(1167, 2).Value = "MARATHON"
$ws.Cells.Item(1167, 7).Value = 579.05

$ws.Cells.Item(1168, 2).Value = "MOL"
$ws.Cells.Item(1168, 7).Value = 84.05

$ws.Cells.Item(1169, 2).Value = "NAGAFERT"
$ws.Cells.Item(1169, 7).Value = 12.94

$ws.Cells.Item(1170, 2).Value = "NFL"
$ws.Cells.Item(1170, 7).Value = 127.11

$ws.Cells.Item(1171, 2).Value = "ORIENTHOT"
$ws.Cells.Item(1171, 7).Value = 139.35

$ws.Cells.Item(1172, 2).Value = "PARADEEP"
$ws.Cells.Item(1172, 7).Value = 82.75

$ws.Cells.Item(1173, 2).Value = "PNC"
$ws.Cells.Item(1173, 7).Value = 66.36

$ws.Cells.Item(1174, 2).Value = "PRSMJOHNSN"
$ws.Cells.Item(1174, 7).Value = 164.75

$ws.Cells.Item(1175, 2).Value = "QUICKHEAL"
$ws.Cells.Item(1175, 7).Value = 518.85

$ws.Cells.Item(1176, 2).Value = "RAJESHEXPO"
$ws.Cells.Item(1176, 7).Value = 286.3

$ws.Cells.Item(1177, 2).Value = "RALLIS"
$ws.Cells.Item(1177, 7).Value = 314.75

$ws.Cells.Item(1178, 2).Value = "RAMASTEEL"
$ws.Cells.Item(1178, 7).Value = 11.75

$ws.Cells.Item(1179, 2).Value = "RCF"
$ws.Cells.Item(1179, 7).Value = 189.18

$ws.Cells.Item(1180, 2).Value = "RKEC"
$ws.Cells.Item(1180, 7).Value = 105.03

$ws.Cells.Item(1181, 2).Value = "RPOWER"
$ws.Cells.Item(1181, 7).Value = 29.55

$ws.Cells.Item(1182, 2).Value = "SAGCEM"
$ws.Cells.Item(1182, 7).Value = 243.73

$ws.Cells.Item(1183, 1).Value = "25/06/2024"

